$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.005.53"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "3.763.52"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "632.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").Value = "3.761.58"
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").Value = "4.398.19"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "3.766.76"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "69.025.65"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "461.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.705"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("E24").Value = "  -4.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.16%  "
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "3.914.97"
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.174"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +18.48%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "3.718.56"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.101"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("E40").Value = "  +3.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.962"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "157.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("E46").Value = "  +5.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
